$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'69.772.96"
$ws.Range("E2").Value = "  +1.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.928.37"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'528.36"
$ws.Range("E5").Value = "  +8.44%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'145.16"
$ws.Range("E6").Value = "  -0.61%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.18%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.727"
$ws.Range("E9").Value = "  -0.91%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.174"
$ws.Range("E10").Value = "  +3.94%  "

# Row 11 - ShibaInu
$ws.Range("D11").Value = "'0.0000335"
$ws.Range("E11").Value = "  -3.07%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'42.49"
$ws.Range("E12").Value = "  -1.80%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'4.556.07"
$ws.Range("E13").Value = "  +0.36%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'10.34"
$ws.Range("E14").Value = "  -4.13%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "'3.937.40"
$ws.Range("E15").Value = "  +0.64%  "

# Row 16 - was TRON, now Uniswap
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'13.97"
$ws.Range("E16").Value = "  -2.71%  "

# Row 17 - was Uniswap, now TRON
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.136"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +6.83%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'19.85"
$ws.Range("E19").Value = "  -1.06%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "'69.665.73"
$ws.Range("E20").Value = "  +1.83%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'434.34"
$ws.Range("E21").Value = "  -0.06%  "

# Row 22 - ImmutableX
$ws.Range("D22").Value = "'3.37"
$ws.Range("E22").Value = "  -3.41%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "'14.48"
$ws.Range("E23").Value = "  -4.57%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "'4.11"
$ws.Range("E24").Value = "  +13.28%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'88.06"
$ws.Range("E25").Value = "  -0.29%  "

# Row 26 - RenderToken
$ws.Range("D26").Value = "'11.62"
$ws.Range("E26").Value = "  +1.45%  "

# Row 27 - Filecoin
$ws.Range("D27").Value = "'10.73"
$ws.Range("E27").Value = "  -4.20%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'36.58"
$ws.Range("E28").Value = "  -3.81%  "

# Row 29 - Bittensor
$ws.Range("D29").Value = "'698.10"
$ws.Range("E29").Value = "  -3.02%  "

# Row 30 - Cosmos
$ws.Range("D30").Value = "'13.25"
$ws.Range("E30").Value = "  -4.02%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -2.29%  "

# Row 32 - Toncoin
$ws.Range("E32").Value = "  -3.15%  "

# Row 33 - OKB
$ws.Range("D33").Value = "'69.61"
$ws.Range("E33").Value = "  +14.34%  "

# Row 34 - TheGraph
$ws.Range("D34").Value = "'0.455"
$ws.Range("E34").Value = "  +15.81%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -3.87%  "

# Row 36 - InjectiveProtocol
$ws.Range("D36").Value = "'40.20"
$ws.Range("E36").Value = "  -3.00%  "

# Row 37 - PEPE
$ws.Range("D37").Value = "'0.0" + [char]0x2083 + "0841"
$ws.Range("E37").Value = "  -4.93%  "

# Row 38 - Kaspa
$ws.Range("D38").Value = "'0.148"
$ws.Range("E38").Value = "  +1.35%  "

# Row 39 - Dai
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.17%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  -0.04%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "'0.0481"
$ws.Range("E41").Value = "  -2.21%  "

# Row 42 - WEMIXToken
$ws.Range("D42").Value = "'3.10"
$ws.Range("E42").Value = "  +3.63%  "

# Row 43 - Fetch.AI
$ws.Range("D43").Value = "'2.77"
$ws.Range("E43").Value = "  -6.31%  "

# Row 44 - ThetaToken
$ws.Range("D44").Value = "'2.97"
$ws.Range("E44").Value = "  -4.94%  "

# Row 45 - Stacks
$ws.Range("D45").Value = "'3.18"
$ws.Range("E45").Value = "  +12.21%  "

# Row 46 - was ApeXProtocol, now Stellar
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.142"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47 - was Stellar, now ApeXProtocol
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.36"
$ws.Range("E47").Value = "  -1.10%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").Value = "'0.0" + [char]0x2086 + "0349"
$ws.Range("E48").Value = "  +3.28%  "

# Row 49 - LidoDAOToken
$ws.Range("D49").Value = "'3.32"
$ws.Range("E49").Value = "  -2.84%  "

# Row 50 - Monero
$ws.Range("D50").Value = "'144.79"
$ws.Range("E50").Value = "  -0.37%  "

# Row 51 - ARBITRUM
$ws.Range("E51").Value = "  -2.95%  "
